$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.646.46"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "1.598.96"
$ws.Range("E3").Value = "  +0.89%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("E6").Value = "  +1.48%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("E9").Value = "  -1.31%  "
$ws.Range("E10").Value = "  +1.03%  "
$ws.Range("E11").Value = "  +0.38%  "
$ws.Range("D12").Value = "1.822.94"
$ws.Range("E12").Value = "  +0.73%  "
$ws.Range("D13").Value = "1.604.85"
$ws.Range("E13").Value = "  +1.04%  "
$ws.Range("E14").Value = "  -0.49%  "
$ws.Range("E15").Value = "  -1.16%  "
$ws.Range("E16").Value = "  +1.56%  "
$ws.Range("D17").Value = "26.640.67"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("D18").Value = "0.0₃0730"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("E21").Value = "  +1.33%  "
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("E23").Value = "  -2.40%  "
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("E25").Value = "  -0.31%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  -2.99%  "
$ws.Range("E28").Value = "  +2.47%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("E30").Value = "  +0.97%  "
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("E32").Value = "  -0.20%  "
$ws.Range("E33").Value = "  +1.38%  "
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("D35").Value = "1.295.49"
$ws.Range("E35").Value = "  -0.66%  "
$ws.Range("E36").Value = "  +0.65%  "
$ws.Range("E37").Value = "  -1.09%  "
$ws.Range("E38").Value = "  -0.53%  "
$ws.Range("E39").Value = "  +3.20%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("E41").Value = "  +2.46%  "
$ws.Range("E42").Value = "  +1.84%  "
$ws.Range("E43").Value = "  +0.50%  "
$ws.Range("E44").Value = "  +1.69%  "
$ws.Range("D45").Value = "1.735.52"
$ws.Range("E45").Value = "  +0.65%  "
$ws.Range("E46").Value = "  +7.21%  "
$ws.Range("E47").Value = "  +1.51%  "
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("E49").Value = "  +2.81%  "
$ws.Range("E50").Value = "  -0.46%  "
$ws.Range("E51").Value = "  +0.89%  "

# Force text storage for numeric-looking price values (avoid Excel auto number conversion)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.512"
$ws.Range("D6").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.66"
$ws.Range("D10").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.84"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "209.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.662"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.93"
$ws.Range("D34").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.49"
$ws.Range("D37").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.788"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.84"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.890"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.19"
$ws.Range("D47").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.52"
$ws.Range("D51").Style = "Normal"
